# Weekly price update: a new Ajo (garlic) quote for
# "Terminal Hortofrutícola Agro Chillán" is inserted as a new row right
# after the existing row for 2021-12-06 (serial 44536), pushing every
# subsequent record down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 154; Excel shifts rows 154..203 down to
# 155..204 and carries the row-154 formatting (e.g. the date style on
# column D) onto the freshly inserted row.
$ws.Rows.Item(154).Insert()

# Populate the new row with the new weekly observation.
$ws.Cells.Item(154, 1).Value  = 7
$ws.Cells.Item(154, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(154, 3).Value  = "Ñuble"
$ws.Cells.Item(154, 4).Value  = 44627
$ws.Cells.Item(154, 5).Value  = 16
$ws.Cells.Item(154, 6).Value  = 100112003
$ws.Cells.Item(154, 7).Value  = "Ajo"
$ws.Cells.Item(154, 8).Value  = "Chino"
$ws.Cells.Item(154, 9).Value  = "Primera"
$ws.Cells.Item(154, 10).Value = 60
$ws.Cells.Item(154, 11).Value = 19000
$ws.Cells.Item(154, 12).Value = 20000
$ws.Cells.Item(154, 13).Value = 19500
$ws.Cells.Item(154, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(154, 15).Value = "China"
$ws.Cells.Item(154, 16).Value = 1950
$ws.Cells.Item(154, 17).Value = 10
$ws.Cells.Item(154, 18).Value = "Hortaliza"
